# Weekly crime report update: new week of data collected
# (commit message: "New crime data collected")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title block: Volume/Number and date range ---
# "Volume 30   Number  35" -> "Volume 30   Number  36"
$ws.Range("A8").Value = "Volume 30   Number  36"
# "Report Covering the Week  8/28/2023  Through  9/3/2023"
#   -> "Report Covering the Week  9/4/2023  Through  9/10/2023"
$ws.Range("C9").Value = "Report Covering the Week  9/4/2023  Through  9/10/2023"

# --- Crime Complaints table (rows 14-30): Week to Date / 28 Day / Year to Date figures ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = -66.666666666666
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 49
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = -2
$ws.Range("L14").Value = -25.757575757575
$ws.Range("M14").Value = -49.484536082474
$ws.Range("N14").Value = -85.588235294117

$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = -60
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = -47.826086956521
$ws.Range("I15").Value = 154
$ws.Range("J15").Value = 171
$ws.Range("K15").Value = -9.941520467836
$ws.Range("L15").Value = 0.653594771241
$ws.Range("M15").Value = 1.986754966887
$ws.Range("N15").Value = -64.516129032258

$ws.Range("C16").Value = 55
$ws.Range("D16").Value = 62
$ws.Range("E16").Value = -11.290322580645
$ws.Range("F16").Value = 203
$ws.Range("G16").Value = 225
$ws.Range("H16").Value = -9.777777777777
$ws.Range("I16").Value = 1689
$ws.Range("J16").Value = 1817
$ws.Range("K16").Value = -7.044578976334
$ws.Range("L16").Value = 25.111111111111
$ws.Range("M16").Value = -29.093198992443
$ws.Range("N16").Value = -85.129424194400

$ws.Range("C17").Value = 94
$ws.Range("D17").Value = 104
$ws.Range("E17").Value = -9.615384615384
$ws.Range("F17").Value = 306
$ws.Range("G17").Value = 354
$ws.Range("H17").Value = -13.559322033898
$ws.Range("I17").Value = 2951
$ws.Range("J17").Value = 2931
$ws.Range("K17").Value = 0.682360968952
$ws.Range("L17").Value = 21.590440873506
$ws.Range("M17").Value = 25.201527365294
$ws.Range("N17").Value = -51.207010582010

$ws.Range("C18").Value = 28
$ws.Range("D18").Value = 38
$ws.Range("E18").Value = -26.315789473684
$ws.Range("F18").Value = 142
$ws.Range("G18").Value = 173
$ws.Range("H18").Value = -17.919075144508
$ws.Range("I18").Value = 1419
$ws.Range("J18").Value = 1656
$ws.Range("K18").Value = -14.311594202898
$ws.Range("L18").Value = 7.826747720364
$ws.Range("M18").Value = -34.030683403068
$ws.Range("N18").Value = -83.066825775656

$ws.Range("C19").Value = 104
$ws.Range("D19").Value = 119
$ws.Range("E19").Value = -12.605042016806
$ws.Range("F19").Value = 434
$ws.Range("G19").Value = 492
$ws.Range("H19").Value = -11.788617886178
$ws.Range("I19").Value = 4009
$ws.Range("J19").Value = 4103
$ws.Range("K19").Value = -2.291006580550
$ws.Range("L19").Value = 30.671447196870
$ws.Range("M19").Value = 40.223854494578
$ws.Range("N19").Value = -15.171392297926

$ws.Range("C20").Value = 29
$ws.Range("D20").Value = 37
$ws.Range("E20").Value = -21.621621621621
$ws.Range("F20").Value = 139
$ws.Range("G20").Value = 166
$ws.Range("H20").Value = -16.265060240963
$ws.Range("I20").Value = 1245
$ws.Range("J20").Value = 1253
$ws.Range("K20").Value = -0.638467677573
$ws.Range("L20").Value = 20.406189555125
$ws.Range("M20").Value = 26.782077393075
$ws.Range("N20").Value = -80.810727496917

$ws.Range("C21").Value = 313
$ws.Range("D21").Value = 368
$ws.Range("E21").Value = -14.945652173913
$ws.Range("F21").Value = 1241
$ws.Range("G21").Value = 1438
$ws.Range("H21").Value = -13.699582753824
$ws.Range("I21").Value = 11516
$ws.Range("J21").Value = 11981
$ws.Range("K21").Value = -3.881145146481
$ws.Range("L21").Value = 22.328446993839
$ws.Range("M21").Value = 4.891155842972
$ws.Range("N21").Value = -69.513421930428

$ws.Range("C22").Value = 4
$ws.Range("E22").Value = -42.857142857142
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = -24
$ws.Range("I22").Value = 197
$ws.Range("J22").Value = 244
$ws.Range("K22").Value = -19.262295081967
$ws.Range("L22").Value = 18.674698795180
$ws.Range("M22").Value = -33.220338983050

$ws.Range("C23").Value = 29
$ws.Range("D23").Value = 31
$ws.Range("E23").Value = -6.451612903225
$ws.Range("F23").Value = 107
$ws.Range("G23").Value = 127
$ws.Range("H23").Value = -15.748031496063
$ws.Range("I23").Value = 1095
$ws.Range("J23").Value = 1079
$ws.Range("K23").Value = 1.482854494902
$ws.Range("L23").Value = 11.054766734279
$ws.Range("M23").Value = 36.024844720496

$ws.Range("C24").Value = 254
$ws.Range("D24").Value = 236
$ws.Range("E24").Value = 7.627118644067
$ws.Range("F24").Value = 979
$ws.Range("G24").Value = 1088
$ws.Range("H24").Value = -10.018382352941
$ws.Range("I24").Value = 8687
$ws.Range("J24").Value = 9174
$ws.Range("K24").Value = -5.308480488336
$ws.Range("L24").Value = 26.466734604746
$ws.Range("M24").Value = 22.715072750388

$ws.Range("C25").Value = 136
$ws.Range("D25").Value = 99
$ws.Range("E25").Value = 37.373737373737
$ws.Range("F25").Value = 495
$ws.Range("G25").Value = 424
$ws.Range("H25").Value = 16.745283018867
$ws.Range("I25").Value = 4296
$ws.Range("J25").Value = 4163
$ws.Range("K25").Value = 3.194811434061
$ws.Range("L25").Value = 35.691724573594
$ws.Range("M25").Value = -23.217158176943

$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = -34.482758620689
$ws.Range("I26").Value = 239
$ws.Range("J26").Value = 259
$ws.Range("K26").Value = -7.722007722007
$ws.Range("L26").Value = -8.778625954198

$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = 66.666666666666
$ws.Range("F27").Value = 66
$ws.Range("G27").Value = 44
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 445
$ws.Range("J27").Value = 435
$ws.Range("K27").Value = 2.298850574712
$ws.Range("L27").Value = -7.484407484407

$ws.Range("C28").Value = 11
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = 37.5
$ws.Range("F28").Value = 21
$ws.Range("G28").Value = 34
$ws.Range("H28").Value = -38.235294117647
$ws.Range("I28").Value = 167
$ws.Range("J28").Value = 255
$ws.Range("K28").Value = -34.509803921568
$ws.Range("L28").Value = -45.065789473684
$ws.Range("M28").Value = -55.466666666666
$ws.Range("N28").Value = -87.863372093023

$ws.Range("C29").Value = 5
$ws.Range("E29").Value = -28.571428571428
$ws.Range("G29").Value = 28
$ws.Range("H29").Value = -46.428571428571
$ws.Range("I29").Value = 142
$ws.Range("J29").Value = 213
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -41.078838174273
$ws.Range("M29").Value = -52.823920265780
$ws.Range("N29").Value = -88.557614826752

# Row 30 (Hate Crimes): C30 was the literal string "0" (no complaints recorded);
# this week it becomes a real numeric figure, matching the numeric formatting of
# the rest of the column.
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = -77.777777777777
$ws.Range("I30").Value = 42
$ws.Range("J30").Value = 55
$ws.Range("K30").Value = -23.636363636363
$ws.Range("L30").Value = -4.545454545454
